$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (2-6); the upload template now ships empty,
# header-only, with the stale shared strings swept away as a side effect.
$ws.Rows("2:6").Delete()

# Insert a new column for "Supplier Code" before the existing "Supplier"
# column (old C), pushing Supplier -> D and Defect Types -> E.
$ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "Supplier Code"

# Match the new column widths for the Supplier Code / Supplier columns.
$ws.Columns("C:C").ColumnWidth = 13.3
$ws.Columns("D:D").ColumnWidth = 14.1

# Restore the last active-cell selection recorded for the sheet.
$ws.Range("E7").Select()
